# TC03_Bento_Filter_EndocrineTherapy-OFS.xlsx - "Fixed Bento 80 Test scripts"
#
# The three Neo4j/Cypher queries stored on the "startup" sheet (column B,
# rows 2-4) are updated so each query explicitly sorts its results and caps
# them at 100 rows (the "Bento 80" paging fix). The Files query's existing
# lowercase "order by" is normalized to "order By ... ASC LIMIT 100"; the
# Cases and Samples queries (which previously had no ordering at all) gain
# a new "order By ... ASC LIMIT 100" clause.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CasesTab): query in B2 gets a new trailing ORDER BY / LIMIT clause.
$casesQuery = $ws.Cells.Item(2, 2).Value()
$casesQuery = $casesQuery + "`n order By ss.study_subject_id ASC LIMIT 100"
$ws.Cells.Item(2, 2).Value = $casesQuery

# Row 3 (SamplesTab): query in B3 gets a new trailing ORDER BY / LIMIT clause.
$samplesQuery = $ws.Cells.Item(3, 2).Value()
$samplesQuery = $samplesQuery + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Cells.Item(3, 2).Value = $samplesQuery

# Row 4 (FilesTab): existing "order by f.file_name" is upgraded in place.
$filesQuery = $ws.Cells.Item(4, 2).Value()
$filesQuery = $filesQuery.Replace("order by f.file_name", "order By f.file_name ASC LIMIT 100")
$ws.Cells.Item(4, 2).Value = $filesQuery

# The extra wrapped line of text grows rows 2 & 3 slightly (wrap-text style).
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 360

# Scroll the view down a row and leave a single-cell selection on B2.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("B2").Select()
